# "test P7 with -10 percent"
# Applies the recorded value changes (re-run of the simulation with a
# -10% perturbation) to the results workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "general": summary/objective values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("general")
$ws.Cells.Item(3, 2).Value = 433.41074364326
$ws.Cells.Item(4, 2).Value = 0.01799988746643066
$ws.Cells.Item(6, 2).Value = 33.93074364325894
$ws.Cells.Item(7, 2).Value = 2.580484450641003
$ws.Cells.Item(8, 2).Value = 2.580484450641003
$ws.Cells.Item(9, 2).Value = 235.100000000001
$ws.Cells.Item(10, 2).Value = 164.38

# ---------------------------------------------------------------
# Sheet "alpha": new rows 2-6 (j, s, alpha)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("alpha")
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 2
$ws.Cells.Item(3, 3).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 3
$ws.Cells.Item(4, 3).Value = 1
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 4
$ws.Cells.Item(5, 3).Value = 1
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 2).Value = 5
$ws.Cells.Item(6, 3).Value = 1

# ---------------------------------------------------------------
# Sheet "x": column B (j) re-assignment
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("x")
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(3, 2).Value = 3
$ws.Cells.Item(4, 2).Value = 9
$ws.Cells.Item(5, 2).Value = 7
$ws.Cells.Item(8, 2).Value = 11
$ws.Cells.Item(9, 2).Value = 13
$ws.Cells.Item(10, 2).Value = 5
$ws.Cells.Item(11, 2).Value = 10
$ws.Cells.Item(12, 2).Value = 4
$ws.Cells.Item(13, 2).Value = 8
$ws.Cells.Item(14, 2).Value = 2

# ---------------------------------------------------------------
# Sheet "U": column B (t) tweaks
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("U")
$ws.Cells.Item(2, 2).Value = 1
$ws.Cells.Item(8, 2).Value = 3
$ws.Cells.Item(11, 2).Value = 3

# ---------------------------------------------------------------
# Sheet "TBar": column B (TBar) recomputed values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("TBar")
$ws.Cells.Item(3, 2).Value = 7.38807534940317
$ws.Cells.Item(4, 2).Value = 26.95051889297535
$ws.Cells.Item(6, 2).Value = 24.76592070603971
$ws.Cells.Item(7, 2).Value = 24.16886835983306
$ws.Cells.Item(9, 2).Value = 20
$ws.Cells.Item(10, 2).Value = 22.45367071955468
$ws.Cells.Item(11, 2).Value = 20
$ws.Cells.Item(12, 2).Value = 22.01159140980467
$ws.Cells.Item(13, 2).Value = 28.28184163802894
$ws.Cells.Item(14, 2).Value = 29.53100334361635
$ws.Cells.Item(15, 2).Value = 25.35398438790795

# ---------------------------------------------------------------
# Sheet "y": new rows 2-6 (i, j, s, y)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("y")
$ws.Cells.Item(2, 1).Value = 2
$ws.Cells.Item(2, 2).Value = 12
$ws.Cells.Item(2, 3).Value = 1
$ws.Cells.Item(2, 4).Value = 1
$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 12
$ws.Cells.Item(3, 3).Value = 2
$ws.Cells.Item(3, 4).Value = 1
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = 12
$ws.Cells.Item(4, 3).Value = 3
$ws.Cells.Item(4, 4).Value = 1
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = 12
$ws.Cells.Item(5, 3).Value = 4
$ws.Cells.Item(5, 4).Value = 1
$ws.Cells.Item(6, 1).Value = 2
$ws.Cells.Item(6, 2).Value = 12
$ws.Cells.Item(6, 3).Value = 5
$ws.Cells.Item(6, 4).Value = 1

# ---------------------------------------------------------------
# Sheet "Q": column C (Q) recomputed values, rows 7-71
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("Q")
$ws.Cells.Item(7, 3).Value = 97.1700000000008
$ws.Cells.Item(8, 3).Value = 99.27000000000081
$ws.Cells.Item(9, 3).Value = 100.1150000000008
$ws.Cells.Item(10, 3).Value = 99.0400000000008
$ws.Cells.Item(11, 3).Value = 97.9800000000008
$ws.Cells.Item(12, 3).Value = 319.6700000000007
$ws.Cells.Item(13, 3).Value = 323.35
$ws.Cells.Item(14, 3).Value = 324.5350000000007
$ws.Cells.Item(15, 3).Value = 329.9
$ws.Cells.Item(16, 3).Value = 320.0950000000008
$ws.Cells.Item(21, 3).Value = 39.43499999999941
$ws.Cells.Item(22, 3).Value = 226.0399999999994
$ws.Cells.Item(23, 3).Value = 247.1799999999994
$ws.Cells.Item(24, 3).Value = 221.8549999999994
$ws.Cells.Item(25, 3).Value = 238.4549999999994
$ws.Cells.Item(26, 3).Value = 224.4749999999995
$ws.Cells.Item(27, 3).Value = 224.1799999999995
$ws.Cells.Item(28, 3).Value = 224.6649999999995
$ws.Cells.Item(29, 3).Value = 201.1149999999995
$ws.Cells.Item(30, 3).Value = 218.9699999999995
$ws.Cells.Item(31, 3).Value = 207.1049999999995
$ws.Cells.Item(35, 3).Value = 146.3249999999992
$ws.Cells.Item(37, 3).Value = 180.2450000000007
$ws.Cells.Item(38, 3).Value = 183.9900000000007
$ws.Cells.Item(39, 3).Value = 178.0900000000007
$ws.Cells.Item(40, 3).Value = 188.8100000000007
$ws.Cells.Item(41, 3).Value = 179.8350000000007
$ws.Cells.Item(42, 3).Value = 153.4099999999999
$ws.Cells.Item(43, 3).Value = 167.1249999999999
$ws.Cells.Item(44, 3).Value = 139.5349999999999
$ws.Cells.Item(45, 3).Value = 154.5
$ws.Cells.Item(46, 3).Value = 143.6599999999999
$ws.Cells.Item(47, 3).Value = 85.48500000000051
$ws.Cells.Item(48, 3).Value = 87.9650000000005
$ws.Cells.Item(49, 3).Value = 79.71500000000052
$ws.Cells.Item(50, 3).Value = 90.7300000000005
$ws.Cells.Item(51, 3).Value = 84.73000000000052
$ws.Cells.Item(52, 3).Value = 238.195
$ws.Cells.Item(53, 3).Value = 242.67
$ws.Cells.Item(54, 3).Value = 239.82
$ws.Cells.Item(55, 3).Value = 249.465
$ws.Cells.Item(56, 3).Value = 232.75
$ws.Cells.Item(57, 3).Value = 226.0399999999994
$ws.Cells.Item(58, 3).Value = 247.1799999999994
$ws.Cells.Item(59, 3).Value = 221.8549999999994
$ws.Cells.Item(60, 3).Value = 238.4549999999994
$ws.Cells.Item(61, 3).Value = 224.4749999999995
$ws.Cells.Item(62, 3).Value = 319.6700000000007
$ws.Cells.Item(63, 3).Value = 323.35
$ws.Cells.Item(64, 3).Value = 324.5350000000007
$ws.Cells.Item(65, 3).Value = 329.9
$ws.Cells.Item(66, 3).Value = 320.0950000000008
$ws.Cells.Item(67, 3).Value = 224.1799999999995
$ws.Cells.Item(68, 3).Value = 224.6649999999995
$ws.Cells.Item(69, 3).Value = 201.1149999999995
$ws.Cells.Item(70, 3).Value = 218.9699999999995
$ws.Cells.Item(71, 3).Value = 207.1049999999995

# ---------------------------------------------------------------
# Sheet "R": column C (R) recomputed values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("R")
$ws.Cells.Item(3, 3).Value = 0
$ws.Cells.Item(5, 3).Value = 0
$ws.Cells.Item(7, 3).Value = 19.67
$ws.Cells.Item(8, 3).Value = 23.35
$ws.Cells.Item(9, 3).Value = 24.53500000000051
$ws.Cells.Item(10, 3).Value = 29.9
$ws.Cells.Item(11, 3).Value = 20.095

# ---------------------------------------------------------------
# Sheet "L": column C (L) recomputed values
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("L")
$ws.Cells.Item(2, 3).Value = 12.775
$ws.Cells.Item(3, 3).Value = 18.32
$ws.Cells.Item(4, 3).Value = 13.155
$ws.Cells.Item(5, 3).Value = 20.115
$ws.Cells.Item(6, 3).Value = 17.825
$ws.Cells.Item(32, 3).Value = 0
$ws.Cells.Item(33, 3).Value = 0
$ws.Cells.Item(34, 3).Value = 0
$ws.Cells.Item(35, 3).Value = 0
$ws.Cells.Item(36, 3).Value = 0
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(51, 3).Value = 0

# ---------------------------------------------------------------
# Sheet "rho": drop rows 7-8 (i=5 rows no longer present)
# ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("rho")
$ws.Range("A7:C8").Delete()
